# Edit corresponding to commit "Tue, Jul 07, 2020 10:05:06 PM":
#  1. Slide 16's summary table switches from the custom "Table_0" style
#     ({1BF8C62C-EA4C-4024-A846-8CF62AFC9379}) to PowerPoint's built-in
#     "Medium Style 2 - Accent 1" table style ({C7A7ABAF-AC3E-4720-988E-E961B28B7A5B}).
#  2. The deck's theme colour palette is switched from the custom "Integral"
#     palette over to the stock Office theme palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{C7A7ABAF-AC3E-4720-988E-E961B28B7A5B}")
    }
}

# --- 2. Theme colours: Integral -> Office -----------------------------------
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

$officeColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = $officeColors[$i - 1]
}
